# "update scripts wuth new tpm"
# Refresh the ligand/receptor expression + specificity metrics (columns G-J,
# M-T) in the Angpt2-Tek LR-pair sheet with the newly recomputed TPM-based
# values. Columns A-F and K-L (cluster labels / expressing-cell counts) are
# unaffected by the TPM recalculation and are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.107822000000001
$ws.Range("H2").Value = 24.323466
$ws.Range("I2").Value = 0.4676336537051783
$ws.Range("J2").Value = 0.4676336537051783
$ws.Range("M2").Value = 53.457377
$ws.Range("N2").Value = 160.372131
$ws.Range("O2").Value = 0.6217639481372091
$ws.Range("P2").Value = 0.6217639481372091
$ws.Range("Q2").Value = 433.4228973028941
$ws.Range("R2").Value = 3900.806075726046
$ws.Range("S2").Value = 0.2907577468095601
$ws.Range("T2").Value = 0.2907577468095601
$ws.Range("G3").Value = 8.107822000000001
$ws.Range("H3").Value = 24.323466
$ws.Range("I3").Value = 0.4676336537051783
$ws.Range("J3").Value = 0.4676336537051783
$ws.Range("O3").Value = 0.3615335470438062
$ws.Range("P3").Value = 0.3615335470438062
$ws.Range("Q3").Value = 252.019947282854
$ws.Range("R3").Value = 2268.179525545686
$ws.Range("S3").Value = 0.1690652535410881
$ws.Range("T3").Value = 0.1690652535410881
$ws.Range("G4").Value = 8.107822000000001
$ws.Range("H4").Value = 24.323466
$ws.Range("I4").Value = 0.4676336537051783
$ws.Range("J4").Value = 0.4676336537051783
$ws.Range("O4").Value = 0.01670250481898457
$ws.Range("P4").Value = 0.01670250481898457
$ws.Range("Q4").Value = 11.64308103187467
$ws.Range("R4").Value = 104.787729286872
$ws.Range("S4").Value = 0.007810653354530105
$ws.Range("T4").Value = 0.007810653354530104
$ws.Range("I5").Value = 0.05150839004025344
$ws.Range("J5").Value = 0.05150839004025344
$ws.Range("M5").Value = 53.457377
$ws.Range("N5").Value = 160.372131
$ws.Range("O5").Value = 0.6217639481372091
$ws.Range("P5").Value = 0.6217639481372091
$ws.Range("Q5").Value = 47.74018180635266
$ws.Range("R5").Value = 429.6616362571739
$ws.Range("S5").Value = 0.03202605995361928
$ws.Range("T5").Value = 0.03202605995361928
$ws.Range("I6").Value = 0.05150839004025344
$ws.Range("J6").Value = 0.05150839004025344
$ws.Range("O6").Value = 0.3615335470438062
$ws.Range("P6").Value = 0.3615335470438062
$ws.Range("S6").Value = 0.01862201095376869
$ws.Range("T6").Value = 0.01862201095376869
$ws.Range("I7").Value = 0.05150839004025344
$ws.Range("J7").Value = 0.05150839004025344
$ws.Range("O7").Value = 0.01670250481898457
$ws.Range("P7").Value = 0.01670250481898457
$ws.Range("S7").Value = 0.0008603191328654702
$ws.Range("T7").Value = 0.0008603191328654702
$ws.Range("I8").Value = 0.4808579562545683
$ws.Range("J8").Value = 0.4808579562545682
$ws.Range("M8").Value = 53.457377
$ws.Range("N8").Value = 160.372131
$ws.Range("O8").Value = 0.6217639481372091
$ws.Range("P8").Value = 0.6217639481372091
$ws.Range("Q8").Value = 445.6797472544594
$ws.Range("R8").Value = 4011.117725290134
$ws.Range("S8").Value = 0.2989801413740298
$ws.Range("T8").Value = 0.2989801413740297
$ws.Range("I9").Value = 0.4808579562545683
$ws.Range("J9").Value = 0.4808579562545682
$ws.Range("O9").Value = 0.3615335470438062
$ws.Range("P9").Value = 0.3615335470438062
$ws.Range("S9").Value = 0.1738462825489495
$ws.Range("T9").Value = 0.1738462825489495
$ws.Range("I10").Value = 0.4808579562545683
$ws.Range("J10").Value = 0.4808579562545682
$ws.Range("O10").Value = 0.01670250481898457
$ws.Range("P10").Value = 0.01670250481898457
$ws.Range("S10").Value = 0.008031532331589
$ws.Range("T10").Value = 0.008031532331589
